# Apply "Doing Updates for Financials" edits to the ASKDF yearly financials sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASKDF")

# --- Income Statement ---
# Total Operating Expenses (row 17)
$ws.Range("D17").Value = 600
$ws.Range("J17").Value = 3700

# Operating Income or Loss (row 18)
$ws.Range("J18").Value = -3700

# Income Before Tax (row 23)
$ws.Range("D23").Value = -600
$ws.Range("J23").Value = -3600

# Income After Tax (row 26)
$ws.Range("D26").Value = -600
$ws.Range("J26").Value = -3600

# Net Income From Continuing Ops (row 27)
$ws.Range("D27").Value = -600
$ws.Range("J27").Value = -3600

# Net Income (row 33)
$ws.Range("D33").Value = -600
$ws.Range("I33").Value = -1300
$ws.Range("J33").Value = -3600

# Net Income Applicable To Common Shares (row 35)
$ws.Range("D35").Value = -600
$ws.Range("I35").Value = -1300
$ws.Range("J35").Value = -3600

# --- Balance Sheet ---
# Total Assets (row 54)
$ws.Range("E54").Value = 400

# Retained Earnings (row 72)
$ws.Range("D72").Value = -26100
$ws.Range("E72").Value = -25500
$ws.Range("F72").Value = -25400
$ws.Range("G72").Value = -25400
$ws.Range("H72").Value = -25400
$ws.Range("I72").Value = -24800
$ws.Range("J72").Value = -23400

# --- Cash Flow Statement ---
# Net Income (row 81)
$ws.Range("D81").Value = -600
$ws.Range("I81").Value = -1300
$ws.Range("J81").Value = -3600

# Total Cash Flow From Operating Activities (row 89)
$ws.Range("E89").Value = -300
$ws.Range("H89").Value = -900
$ws.Range("J89").Value = -3400

# Total Cash Flows From Investing Activities (row 94)
$ws.Range("J94").Value = 3100

$wb.Save()
